$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new "line" entries are inserted ahead of the "extr" series, so every
# row from 8 downward is relabelled one notch further down that series
# (line7/line8, then extr1..extr8) while the numeric measurements shuffle.
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("E9").Value = $true

$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("E12").Value = $true

$ws.Range("B13").Value = "extr4"
$ws.Range("D13").Value = 8

$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11

# New row 16 - re-uses the style of the preceding "A" column cell.
$ws.Range("A16").Value = 14
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

# New row 17
$ws.Range("A17").Value = 15
$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
